$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

$ws.Range("B3").Value = " H R"
$ws.Range("B4").Value = " H R"
$ws.Range("B5").Value = " H R"

$ws.Range("B11").Select()
